$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.217.98'
$ws.Range("E2").Value = '  +0.16%  '

$ws.Range("D3").Value = '2.832.99'
$ws.Range("E3").Value = '  +1.39%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '362.27'
$ws.Range("E5").Value = '  +3.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '112.04'
$ws.Range("E6").Value = '  -2.89%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.569'
$ws.Range("E7").Value = '  +3.83%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.603'
$ws.Range("E9").Value = '  +2.54%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.94'
$ws.Range("E10").Value = '  -3.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0866'
$ws.Range("E11").Value = '  +1.17%  '

$ws.Range("E12").Value = '  +1.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.99'
$ws.Range("E13").Value = '  +0.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.80'
$ws.Range("E14").Value = '  -0.55%  '

$ws.Range("D15").Value = '3.288.71'
$ws.Range("E15").Value = '  +1.58%  '

$ws.Range("D16").Value = '2.829.05'
$ws.Range("E16").Value = '  +1.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.932'
$ws.Range("E17").Value = '  +5.00%  '

$ws.Range("D18").Value = '52.100.49'
$ws.Range("E18").Value = '  -0.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.49'
$ws.Range("E19").Value = '  +3.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.14'
$ws.Range("E20").Value = '  -0.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.41'
$ws.Range("E21").Value = '  +0.47%  '

$ws.Range("D22").Value = '0.0₃0997'
$ws.Range("E22").Value = '  +1.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '272.59'
$ws.Range("E23").Value = '  +1.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.45'
$ws.Range("E24").Value = '  +0.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.83'
$ws.Range("E25").Value = '  +2.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.90'
$ws.Range("E26").Value = '  +0.66%  '

$ws.Range("E27").Value = '  +0.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.32'
$ws.Range("E28").Value = '  +0.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.25'
$ws.Range("E29").Value = '  -0.31%  '

$ws.Range("E30").Value = '  +1.68%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0478'
$ws.Range("E31").Value = '  +2.68%  '

$ws.Range("B32").Value = 'OKB'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '52.32'
$ws.Range("E32").Value = '  +4.29%  '

$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.07'
$ws.Range("E33").Value = '  +2.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.88'
$ws.Range("E34").Value = '  +2.25%  '

$ws.Range("E35").Value = '  +13.25%  '

$ws.Range("E36").Value = '  +2.01%  '

$ws.Range("E37").Value = '  +0.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.27'
$ws.Range("E38").Value = '  +2.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.04'
$ws.Range("E39").Value = '  -2.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.37'
$ws.Range("E40").Value = '  -1.72%  '

$ws.Range("E41").Value = '  +1.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.54'
$ws.Range("E42").Value = '  -2.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '124.68'
$ws.Range("E43").Value = '  -1.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.89'
$ws.Range("E44").Value = '  -1.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.27'
$ws.Range("E45").Value = '  -1.19%  '

$ws.Range("D46").Value = '2.084.53'
$ws.Range("E46").Value = '  +1.52%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.33'
$ws.Range("E47").Value = '  +0.69%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.32'
$ws.Range("E48").Value = '  +2.65%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.88'
$ws.Range("E49").Value = '  +5.74%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.967'
$ws.Range("E50").Value = '  +2.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.20'
$ws.Range("E51").Value = '  +2.69%  '
